$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day's row of profit data (run on 2025-12-11).
# Column A holds the date as plain text (matching the existing rows), so we
# temporarily force a Text number format before writing it, then restore the
# cell's style back to Normal so no stray formatting is left behind.
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "12/11/2025"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = 12749.24
$ws.Range("C17").Value = 0.1957653634455948
$ws.Range("D17").Value = 0.8042346365544052
$ws.Range("E17").Value = -115.03
$ws.Range("F17").Value = -26.44
$ws.Range("G17").Value = -20238.37
$ws.Range("H17").Value = -66.37
$ws.Range("I17").Value = -406.99
$ws.Range("J17").Value = -14.02
